$p = $ppt.ActivePresentation

# Slides 3, 4, and 5 each have an empty "Title 1" placeholder (id=2, type="title")
# that should be removed, as per the commit adding Day 10 - Swagger course
# materials content in its place.
$slideIndexes = @(3, 4, 5)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    $shape = $s.Shapes.Item("Title 1")
    $shape.Delete()
}
